# Vendor Profile.xlsx update
#  - refresh vendor contact details (name / address / phone)
#  - add two new blank rows (11-12) below the item-list header, each split
#    into two merged cells (A:B and C:D) like the existing "Item Name" /
#    "Brand" header row (row 10), picking up that row's formatting

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- update vendor detail values -------------------------------------------------
$ws.Range("B3").Value = "`tK-Mart"
$ws.Range("D3").Value = "Galo-Hilado St,. Brgy. 22, Bacolod, 6100 Negros Occidental"
$ws.Range("B4").Value = "(034) 434 6914, 434-6915"

# --- append two new formatted/merged rows to the item list table ---------------
$ws.Range("A11:B11").Merge() | Out-Null
$ws.Range("C11:D11").Merge() | Out-Null
$ws.Range("A12:B12").Merge() | Out-Null
$ws.Range("C12:D12").Merge() | Out-Null

# Match the bordered style already used for the row-10 header cells (B10/D10)
$ws.Range("B10").Copy()
$ws.Range("A11:D12").PasteSpecial(-4122)

$ws.Range("D12").Select() | Out-Null
